# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update "Conversión del día" note with new Binance rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$nuevoTexto = @"
Conversión del día 💰
✅ Dólar paralelo: 68

Binance
✅ 1000 Bs = 10.0 = 41030.0 pesos
✅ 41030.0 pesos = 9.93 = 933.99 Bs

Promedio competencia
✅ Tasa pesos: 20
✅ Tasa Bs: 20
✅ % Ganancia: 20%
"@

$ws1.Range("A1").Value = $nuevoTexto

# --- tasas: update Binance/transfi reference rates ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 100
$ws2.Range("O10").Value = 4103
$ws2.Range("N12").Value = 4129.99
$ws2.Range("O12").Value = 94.01300000000001
